$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark that currently sits right after
#    the "Properties" heading text.
# ---------------------------------------------------------------------------
$d.Bookmarks.ShowHidden = $true
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) "Type: string " -> "Type: string"  (drop the trailing space)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Type: string ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Type: string", 2)

# ---------------------------------------------------------------------------
# Locate the paragraph that now reads exactly "Type: string" so the two new
# paragraphs land right after it.
# ---------------------------------------------------------------------------
$typeParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd() -eq "Type: string") {
        $typeParaIndex = $i
    }
}

$typePara = $d.Paragraphs($typeParaIndex)

# ---------------------------------------------------------------------------
# 3) New "Channels" heading paragraph (style "berschrift3" == Heading 3).
# ---------------------------------------------------------------------------
$null = $typePara.Range.InsertParagraphAfter()
$headingPara = $d.Paragraphs($typeParaIndex + 1)
$headingPara.Style = "berschrift3"
$headingPara.Range.Text = "Channels"

# ---------------------------------------------------------------------------
# 4) New body paragraph explaining channels.
# ---------------------------------------------------------------------------
$null = $headingPara.Range.InsertParagraphAfter()
$bodyPara = $d.Paragraphs($typeParaIndex + 2)
$bodyPara.Style = "Standard"
$bodyPara.Range.Text = "A channel usually holds more than just one value. That" + [char]8217 + "s why the "

$pos1 = $bodyPara.Range.End - 1
$d.Range($pos1, $pos1).InsertAfter("XmlRpc")

$pos2 = $bodyPara.Range.End - 1
$d.Range($pos2, $pos2).InsertAfter(" method returns an object of type object.")

# Insert the trailing single space FIRST, then drop the empty "_GoBack"
# bookmark right in front of it -- adding a zero-width bookmark exactly at
# the paragraph's trailing edge mis-resolves to the document start, so the
# space run has to exist before the bookmark is added.
$pos3 = $bodyPara.Range.End - 1
$d.Range($pos3, $pos3).InsertAfter(" ")

$bmPos = $pos3
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos))
